# "Início Validação por Horário"
# - Move the DADOS (table "ARRAY") block from B8:C14 up-and-left to A2:B8
# - Re-point the ListObject / table over the new range
# - Add a new "ASCII" worksheet with a binary/ASCII lookup grid
# - Restore view selections for ARRAY / DADOS sheets
# - Make the new ASCII sheet the active tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) DADOS sheet: shift the whole block up 6 rows and left 1 column
# ---------------------------------------------------------------------
$wsDados = $wb.Worksheets.Item("DADOS")
$wsDados.Activate()

$wsDados.Rows("1:6").Delete()
$wsDados.Columns("A").Delete()

# Re-create the small bordered/centred strip that used to sit in the
# (now deleted) decorative row 6 - it now lives alongside the IND=3 row.
$wsDados.Range("C6:G6").HorizontalAlignment = -4108   # xlCenter
$wsDados.Range("C6:G6").VerticalAlignment = -4108     # xlCenter

# The table ("ListObject") needs to be told its range moved too - row
# delete updates it automatically, column delete does not.
$loArray = $wsDados.ListObjects.Item(1)
$loArray.Resize($wsDados.Range("A2:B8"))

$wsDados.Range("C5").Select()

# ---------------------------------------------------------------------
# 2) New "ASCII" worksheet, placed after DADOS
# ---------------------------------------------------------------------
$wsAscii = $wb.Worksheets.Add([System.Type]::Missing, $wsDados)
$wsAscii.Name = "ASCII"

function Set-Row($ws, $row, $colStart, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $colStart + $i).Value = $values[$i]
    }
}

# column bit-weights header row (A:D then F:I -> column E stays blank)
Set-Row $wsAscii 3 1 @(1, 2, 4, 8)
Set-Row $wsAscii 3 6 @(16, 32, 64, 128)

# "H" = 0100 1000 split across the two nibbles
Set-Row $wsAscii 4 1 @(1, 0, 0, 0)
Set-Row $wsAscii 4 6 @(0, 0, 1, 0)

Set-Row $wsAscii 5 1 @(0, 1, 0, 0)
Set-Row $wsAscii 5 6 @(0, 0, 1, 0)

# reversed bit-weights header row for the second mini-table
Set-Row $wsAscii 8 1 @(128, 64, 32, 16)
Set-Row $wsAscii 8 6 @(8, 4, 2, 1)

# "A" = 0100 0001 -> decimal 65
Set-Row $wsAscii 9 1 @(0, 1, 0, 0)
Set-Row $wsAscii 9 6 @(0, 0, 0, 1)
$wsAscii.Cells.Item(9, 10).Value = 65
$wsAscii.Cells.Item(9, 11).Value = "A"

# "B" = 0100 0010 -> decimal 66
Set-Row $wsAscii 10 1 @(0, 1, 0, 0)
Set-Row $wsAscii 10 6 @(0, 0, 1, 0)
$wsAscii.Cells.Item(10, 10).Value = 66
$wsAscii.Cells.Item(10, 11).Value = "B"

$wsAscii.Range("H16").Select()

# ---------------------------------------------------------------------
# 3) Restore the ARRAY sheet's remembered selection
# ---------------------------------------------------------------------
$wsArray = $wb.Worksheets.Item("ARRAY")
$wsArray.Activate()
$wsArray.Range("AD13").Select()

# ---------------------------------------------------------------------
# 4) Leave ASCII as the active sheet/tab
# ---------------------------------------------------------------------
$wsAscii.Activate()
